# Update countries & provincias Spain
# Applies the 25-Jul-2020 14:08 data refresh to the "Pais" sheet:
#   - Updates the "Datos actualizados..." timestamp string.
#   - Refreshes case counts for a number of countries.
#   - A handful of countries swapped rank (the underlying data is sorted
#     by total cases) so their row labels (column A) trade places while
#     the row's statistics are refreshed with the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 14:08"

# --- Countries whose row order / label did not change ------------------
# Row 6 - India
$ws.Range("B6").Value = 1342166
$ws.Range("C6").Value = 5144
$ws.Range("D6").Value = 852497
$ws.Range("E6").Value = 458205
$ws.Range("G6").Value = 58
$ws.Range("H6").Value = 31464

# Row 35 - Bielorrusia
$ws.Range("B35").Value = 67002
$ws.Range("C35").Value = 156
$ws.Range("D35").Value = 60092
$ws.Range("E35").Value = 6380
$ws.Range("G35").Value = 6
$ws.Range("H35").Value = 530

# Row 39 - Kuwait
$ws.Range("B39").Value = 63309
$ws.Range("C39").Value = 684
$ws.Range("D39").Value = 53607
$ws.Range("E39").Value = 9273
$ws.Range("G39").Value = 4
$ws.Range("H39").Value = 429

# Row 42 - Emiratos Arabes Unidos
$ws.Range("B42").Value = 58562
$ws.Range("C42").Value = 313
$ws.Range("D42").Value = 51628
$ws.Range("E42").Value = 6591

# Row 55 - Suiza
$ws.Range("B55").Value = 34302
$ws.Range("C55").Value = 148
$ws.Range("E55").Value = 1825

# Row 67 - Nepal
$ws.Range("B67").Value = 18483
$ws.Range("C67").Value = 109
$ws.Range("D67").Value = 13053
$ws.Range("E67").Value = 5385
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 45

# Row 84 - Senegal
$ws.Range("B84").Value = 9552
$ws.Range("C84").Value = 130
$ws.Range("D84").Value = 6364
$ws.Range("E84").Value = 3001
$ws.Range("G84").Value = 5
$ws.Range("H84").Value = 187

# Row 162 - Vietnam
$ws.Range("B162").Value = 417
$ws.Range("C162").Value = 4
$ws.Range("E162").Value = 52

# Row 179 - Gibraltar
$ws.Range("B179").Value = 185
$ws.Range("C179").Value = 1
$ws.Range("E179").Value = 5

# --- Countries that swapped rank (labels + stats trade rows) -----------
# Rows 80/81: Bulgaria <-> Estado de Palestina
$ws.Range("A80").Value = "Estado de Palestina"
$ws.Range("B80").Value = 10306
$ws.Range("C80").Value = 213
$ws.Range("D80").Value = 3282
$ws.Range("E80").Value = 6953
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 71

$ws.Range("A81").Value = "Bulgaria"
$ws.Range("B81").Value = 10123
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 5252
$ws.Range("E81").Value = 4534
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 337

# Rows 87/88: Consejo Danes para los Refugiados <-> Madagascar
$ws.Range("A87").Value = "Madagascar"
$ws.Range("B87").Value = 8866
$ws.Range("C87").Value = 125
$ws.Range("D87").Value = 5579
$ws.Range("E87").Value = 3209
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 78

$ws.Range("A88").Value = "Consejo Danes para los Refugiados"
$ws.Range("B88").Value = 8801
$ws.Range("C88").Value = 34
$ws.Range("D88").Value = 5305
$ws.Range("E88").Value = 3292
$ws.Range("G88").Value = 3
$ws.Range("H88").Value = 204

# Rows 187/188: Butan <-> Islas Turcas y Caicos
$ws.Range("A187").Value = "Islas Turcas y Caicos"
$ws.Range("B187").Value = 92
$ws.Range("C187").Value = 2
$ws.Range("D187").Value = 28
$ws.Range("E187").Value = 62
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 2

$ws.Range("A188").Value = "Butan"
$ws.Range("B188").Value = 92
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 85
$ws.Range("E188").Value = 7
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0

# Rows 210/211: Groenlandia <-> Islas Malvinas (stats identical, labels only)
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"
